$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Ativação:" date from 01/01/2018 to 01/01/2021.
#    The source text is a plain label (not a real date), so force the
#    cells to Text format first - otherwise Excel auto-converts a
#    dd/mm/yyyy-looking string into a date serial number.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "01/01/2021"

# 2. Update Syllabus row (row 17) with the full English syllabus text
#    (previously it duplicated the short English syllabus text)
$syllabus = "1. Introduction to Operational Research 1.1. Concepts of Operational Research; 1.2. Modeling; 1.3. Structure of Mathematical Models; 1.4. Mathematical techniques in Operational Research; 1.2. Phases of a Study in Operational Research 2. Linear Programming 2.1. Definition 2.2. Formulation of Models 2.3. Graphic Resolution; 3. Simplex method 3.1. Development of the Simplex Method; 3.2. Simplex Method Procedure; 4. Introduction to Graphs and Network Optimization 4.1. Basic Concepts in Graph Theory 4.2. Maximum Flow Problems; 4.3. Minimum Path Problems 5. Case Studies in Linear Programming 5.1. Simple Transport Model 5.2. Model of Designation. 6. Introduction to Queuing Theory 6.1. Queuing Theory Concepts 6.2. Markovian Models"
$ws.Range("B17").Value = $syllabus
$ws.Range("C17").Value = $syllabus

# 3. Update "Método:" row (row 19)
$metodo = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# 4. Update "Critério:" row (row 20)
$criterio = "NF≥ 5,0."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# 5. Update "Norma de recuperação:" row (row 21)
$recuperacao = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao
